# Generate Report for Handback
#
# - Status text "Ready for handoff" -> "Handback transform failed" for the
#   062bac2a-... row (Overview!E3/F3, zh-cn!C3, de-de!C3 all share that text).
# - Populate the (previously empty) "Error Detail" column (P) for that same
#   row on the zh-cn and de-de sheets with the handback/handoff filename
#   mismatch message (language-specific suffix).
# - Widen the "Error Detail" column on zh-cn / de-de so the long message is
#   readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "Handback transform failed"

# Update every cell that shows the old status text for the
# 062bac2a-cd3c-4b7a-a8e4-1b17efc99c81 row.
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Widen the Error Detail column (P) on both locale sheets.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17

# Fill in the Error Detail message explaining the handback/handoff file
# name mismatch, per locale.
$wsZhCn.Range("P3").Value = "Handback file name: itixnz1i.imq is different with handoff file name: 062bac2a-cd3c-4b7a-a8e4-1b17efc99c81.682402cb0ce962f8393324458eeeb31375b66f9a.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: itixnz1i.imq is different with handoff file name: 062bac2a-cd3c-4b7a-a8e4-1b17efc99c81.682402cb0ce962f8393324458eeeb31375b66f9a.de-de."

Write-Host "Handback report generated"
